# Add a new "Sexy Header" entry to the #LocalizedText metadata table, pointing
# at the sexy-script test header that now also exercises the localized text
# table (tables\rococo.tables.test.sxh), per the commit message:
# "Added localized text table to the list of tables scripted with sexy in
#  the carpenter test."

$wb = $excel.ActiveWorkbook

$wsTable = $wb.Worksheets.Item("LocalizedText Table")
$wsMeta  = $wb.Worksheets.Item("#LocalizedText")

# New row appended after the existing "C++ Namespace" row (row 10).
[void]$wsMeta.Activate()
$wsMeta.Range("A11").Value = "Sexy Header"
$wsMeta.Range("B11").Value = "tables\rococo.tables.test.sxh"
[void]$wsMeta.Range("A12").Select()

# Leave the workbook focused back on the LocalizedText Table sheet (the tab
# that was selected before and after the edit), with the cursor moved down
# to C4.
[void]$wsTable.Activate()
[void]$wsTable.Range("C4").Select()
